$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores one weekly pair of rows (Primera/Segunda quality) of
# "Cilantro" price data per iteration, starting at row 36. A new week's
# data (2021-11-25) needs to be inserted at the top of that block (row 36),
# pushing every existing weekly pair down by two rows.
#
# Insert two brand-new blank rows at row 36 (calling Insert twice shifts
# the previously-36/37 pair, and everything below, down by two rows in
# total, growing the sheet from 147 to 149 rows).
$ws.Rows.Item(36).Insert()
$ws.Rows.Item(36).Insert()

# Populate the two newly inserted rows with the same reference data as the
# rest of the "Vega Monumental Concepción - Cilantro" block, dated
# 2021-11-25 (serial 44525).
$ws.Range("A36").Value = 11
$ws.Range("B36").Value = "Vega Monumental Concepción"
$ws.Range("C36").Value = "Bíobío"
$ws.Range("D36").Value = 44525
$ws.Range("E36").Value = 8
$ws.Range("F36").Value = 100112040
$ws.Range("G36").Value = "Cilantro"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 200
$ws.Range("K36").Value = 600
$ws.Range("L36").Value = 700
$ws.Range("M36").Value = 650
$ws.Range("N36").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O36").Value = "Región de Ñuble"
$ws.Range("P36").Value = 650
$ws.Range("Q36").Value = 1
$ws.Range("R36").Value = "Hortaliza"

$ws.Range("A37").Value = 11
$ws.Range("B37").Value = "Vega Monumental Concepción"
$ws.Range("C37").Value = "Bíobío"
$ws.Range("D37").Value = 44525
$ws.Range("E37").Value = 8
$ws.Range("F37").Value = 100112040
$ws.Range("G37").Value = "Cilantro"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Segunda"
$ws.Range("J37").Value = 100
$ws.Range("K37").Value = 500
$ws.Range("L37").Value = 500
$ws.Range("M37").Value = 500
$ws.Range("N37").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O37").Value = "Región de Ñuble"
$ws.Range("P37").Value = 500
$ws.Range("Q37").Value = 1
$ws.Range("R37").Value = "Hortaliza"
